# Re-process the metadata sheet with the newly curated dimensions.
# Rows:
#   1 -> column "variable" identifiers (unchanged)
#   2 -> iaest vocabulary term for each column
#   3 -> role (medida/dim/null)
#   4 -> datatype / URI kind
# Row 5 (mapping-*.xlsx helper row) is removed entirely, and several
# dimension identifiers became "measure" identifiers instead, as per the
# newly curated dimensions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: iaest-measure / iaest-dimension / sdmx-dimension vocabulary ---
$ws.Range("B2").Value = "iaest-measure:nivel-estudios-detalle"
$ws.Range("D2").Value = "sdmx-dimension:refArea"
$ws.Range("E2").Value = "sdmx-dimension:refArea"
$ws.Range("F2").Value = "iaest-measure:nivel-estudios-agregado"
$ws.Range("G2").Value = "sdmx-dimension:refArea"
$ws.Range("I2").Value = "sdmx-dimension:refArea"

# --- Row 3: medida / dim / null role ---
$ws.Range("A3").Value = "medida"
$ws.Range("B3").Value = "medida"
$ws.Range("D3").Value = "dim"
$ws.Range("E3").Value = "dim"
$ws.Range("F3").Value = "medida"
$ws.Range("G3").Value = "dim"
$ws.Range("I3").Value = "dim"

# --- Row 4: datatype / URI reference ---
$ws.Range("A4").Value = "xsd:int"
$ws.Range("B4").Value = "xsd:int"
$ws.Range("D4").Value = "URI-Municipio"
$ws.Range("E4").Value = "URI-Provincia"
$ws.Range("F4").Value = "xsd:int"
$ws.Range("G4").Value = "URI-Comunidad"
$ws.Range("I4").Value = "URI-comarca"

# --- Row 5 is no longer needed (mapping-*.xlsx helper row) ---
$ws.Range("B5").Value = $null
$ws.Range("F5").Value = $null
$ws.Range("G5").Value = $null
$ws.Rows.Item(5).Delete()
